$d = $word.ActiveDocument

$pairs = @(
    @{old="810÷4=202, 2"; new="688÷8=86, 0"},
    @{old="869÷9=96, 5"; new="268÷8=33, 4"},
    @{old="152÷4=38, 0"; new="307÷4=76, 3"},
    @{old="584÷7=83, 3"; new="628÷4=157, 0"},
    @{old="480÷9=53, 3"; new="646÷8=80, 6"},
    @{old="710÷5=142, 0"; new="786÷2=393, 0"},
    @{old="168÷6=28, 0"; new="960÷7=137, 1"},
    @{old="329÷3=109, 2"; new="982÷9=109, 1"},
    @{old="304÷4=76, 0"; new="918÷4=229, 2"},
    @{old="280÷2=140, 0"; new="219÷9=24, 3"},
    @{old="205÷4=51, 1"; new="544÷8=68, 0"},
    @{old="705÷7=100, 5"; new="272÷9=30, 2"},
    @{old="223÷9=24, 7"; new="222÷5=44, 2"},
    @{old="937÷3=312, 1"; new="674÷4=168, 2"},
    @{old="188÷4=47, 0"; new="425÷5=85, 0"},
    @{old="103÷4=25, 3"; new="568÷8=71, 0"},
    @{old="315÷3=105, 0"; new="617÷8=77, 1"},
    @{old="227÷8=28, 3"; new="854÷9=94, 8"},
    @{old="451÷3=150, 1"; new="249÷6=41, 3"},
    @{old="318÷4=79, 2"; new="994÷8=124, 2"},
    @{old="678÷2=339, 0"; new="560÷7=80, 0"},
    @{old="733÷8=91, 5"; new="635÷8=79, 3"},
    @{old="842÷9=93, 5"; new="209÷3=69, 2"},
    @{old="268÷6=44, 4"; new="520÷6=86, 4"},
    @{old="330÷5=66, 0"; new="101÷3=33, 2"}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}

Write-Output "Done replacing $($pairs.Count) values"
